# Commit: "Add files via upload" — re-upload of threshold.xlsx with updated
# threshold values for the "ratio" row (A5="ratio", B5=Min, C5=Max).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 ("ratio"): Min 0.9 -> 0.95, Max 1.5 -> 1.45
$ws.Range("B5").Value = 0.95
$ws.Range("C5").Value = 1.45
